$d = $word.ActiveDocument

# Locate the paragraph that contains the "Ver no Jupiter ..." text; the
# edit removes it together with the blank paragraph right before it and
# the two paragraphs right after it (one blank, one blank-with-page-break)
# that were used to force a page break ahead of the next section.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    # Paragraph just before the "Ver no Jupiter..." paragraph (blank one).
    $startPara = $d.Paragraphs.Item($target - 1)
    # Paragraph two after the "Ver no Jupiter..." paragraph (blank,
    # page-break-before one).
    $endPara = $d.Paragraphs.Item($target + 2)

    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
